{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block that\n// was dropped from the published page, along with the blank paragraph\n// that separated it from the bibliography entry above it.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paras.items;\nfor (let i = 0; i < items.length; i++) {\n  if (targets.indexOf(items[i].text) !== -1) {\n    // Delete this paragraph and the blank paragraph immediately before it\n    // (only once, when hitting the first target paragraph).\n    if (items[i].text === targets[0] && i > 0 && items[i - 1].text === \"\") {\n      items[i - 1].delete();\n    }\n    items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block that was\n# dropped from the published page, along with the blank paragraph that\n# separated it from the bibliography entry above it.\n$d = $word.ActiveDocument\n\nfunction Get-ParaIndexAtPos($doc, $pos) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($pos -ge $p.Range.Start -and $pos -lt $p.Range.End) {\n            return $i\n        }\n    }\n    return $count\n}\n\n# Locate the two paragraphs that must go by their (unique) text.\n$r1 = $d.Content\n$found1 = $r1.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\n\n$r2 = $d.Content\n$found2 = $r2.Find.Execute(\"Contact: luizeleno@usp.br\")\n\nif ($found1 -and $found2) {\n    $idxFirst = Get-ParaIndexAtPos $d $r1.Start\n    $idxLast = Get-ParaIndexAtPos $d $r2.Start\n\n    # Also drop the blank paragraph immediately preceding the block, if any.\n    if ($idxFirst -gt 1) {\n        $prev = $d.Paragraphs.Item($idxFirst - 1)\n        if ($prev.Range.Text.Trim().Length -eq 0) {\n            $idxFirst = $idxFirst - 1\n        }\n    }\n\n    $startPos = $d.Paragraphs.Item($idxFirst).Range.Start\n    $endPos = $d.Paragraphs.Item($idxLast).Range.End\n\n    $delRange = $d.Range($startPos, $endPos)\n    $delRange.Delete()\n}\n"}
